# Modulo de existencias completado
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Encabezado: dia de la factura (13 -> 19)
$ws.Range("E6").Value = "19"

# Numero de factura (A21216407 -> A21216406)
$ws.Range("I3").Value = "A21216406"

# Proveedor (SuperTony Papeleria -> Office Depot)
$ws.Range("H6").Value = "Office Depot"

# Partida / descripcion del primer renglon
$ws.Range("A9").Value = "21201"
$ws.Range("B9").Value = "Tinta de impresora marca Epson"

# Nombre del almacenista (correccion de acento)
$ws.Range("B26").Value = "Tila del Carmen Mendoza Olan"

# Cantidades numericas
$ws.Range("F6").Value = 13
$ws.Range("F9").Value = 5
$ws.Range("H9").Value = 130
